$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$naText = "N/A - Stat tracked as of the 1973-74 ABA Season"

# Row 6 (Brian Taylor, 1972-73 ABA season): stl_per_game / blk_per_game were not tracked yet.
$ws.Range("AQ6").Value = $naText
$ws.Range("AR6").Value = $naText

# Row 7 (Artis Gilmore, 1971-72 ABA season): stl_per_game / blk_per_game were not tracked yet.
# AR7 previously held a stray numeric placeholder value that should be replaced too.
$ws.Range("AQ7").Value = $naText
$ws.Range("AR7").Value = $naText
